# Auto-generated edit script applying scheduled-runner value updates
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 84141
$ws.Range("J3").Value = 84141
$ws.Range("L3").Value = 84141
$ws.Range("N3").Value = -84369
$ws.Range("H53").Value = 588.5
$ws.Range("I53").Value = 544
$ws.Range("J53").Value = 722
$ws.Range("K53").Value = 544
$ws.Range("L53").Value = 722
$ws.Range("M53").Value = 93
$ws.Range("N53").Value = -1996
$ws.Range("H102").Value = 84141
$ws.Range("J102").Value = 84141
$ws.Range("L102").Value = 84141
$ws.Range("N102").Value = -90631
$ws.Range("H129").Value = 2337.0476
$ws.Range("I129").Value = 1916.4706
$ws.Range("J129").Value = 4124.5
$ws.Range("K129").Value = 5749.4118
$ws.Range("L129").Value = 12373.5
$ws.Range("M129").Value = -749.4117999999999
$ws.Range("N129").Value = -22373.5
$ws.Range("H138").Value = 4027.8474
$ws.Range("J138").Value = 5245.943
$ws.Range("L138").Value = 15737.829
$ws.Range("N138").Value = -26017.829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 96079
$ws.Range("J68").Value = 96079
$ws.Range("L68").Value = 96079
$ws.Range("N68").Value = -97701
$ws.Range("H71").Value = 96079
$ws.Range("J71").Value = 96079
$ws.Range("L71").Value = 288237
$ws.Range("N71").Value = -296349
$ws.Range("H81").Value = 98090.375
$ws.Range("J81").Value = 98090.375
$ws.Range("L81").Value = 98090.375
$ws.Range("N81").Value = -100086.375
$ws.Range("H84").Value = 98090.375
$ws.Range("J84").Value = 98090.375
$ws.Range("L84").Value = 294271.125
$ws.Range("N84").Value = -304255.125
$ws.Range("H92").Value = 83625
$ws.Range("J92").Value = 83625
$ws.Range("L92").Value = 83625
$ws.Range("N92").Value = -88617
$ws.Range("H110").Value = 1298.4546
$ws.Range("I110").Value = 1142.5555
$ws.Range("K110").Value = 1142.5555
$ws.Range("M110").Value = 902.4445000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 60610.453
$ws.Range("I88").Value = 29600.5
$ws.Range("J88").Value = 67501.55499999999
$ws.Range("K88").Value = 29600.5
$ws.Range("L88").Value = 67501.55499999999
$ws.Range("M88").Value = -29194.5
$ws.Range("N88").Value = -68313.55499999999
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 60610.453
$ws.Range("I91").Value = 29600.5
$ws.Range("J91").Value = 67501.55499999999
$ws.Range("K91").Value = 29600.5
$ws.Range("L91").Value = 67501.55499999999
$ws.Range("M91").Value = -28196.5
$ws.Range("N91").Value = -70309.55499999999
$ws.Range("H132").Value = 105442
$ws.Range("J132").Value = 105442
$ws.Range("L132").Value = 105442
$ws.Range("N132").Value = -115562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1946.4
$ws.Range("I16").Value = 1722.6
$ws.Range("K16").Value = 1722.6
$ws.Range("M16").Value = -1435.6
$ws.Range("H31").Value = 6337.857
$ws.Range("I31").Value = 5258.885
$ws.Range("J31").Value = 9454.888999999999
$ws.Range("K31").Value = 5258.885
$ws.Range("L31").Value = 9454.888999999999
$ws.Range("M31").Value = -4963.885
$ws.Range("N31").Value = -10044.889
$ws.Range("H34").Value = 6337.857
$ws.Range("I34").Value = 5258.885
$ws.Range("J34").Value = 9454.888999999999
$ws.Range("K34").Value = 5258.885
$ws.Range("L34").Value = 9454.888999999999
$ws.Range("M34").Value = -5056.885
$ws.Range("N34").Value = -9858.888999999999
$ws.Range("H58").Value = 9806092
$ws.Range("I58").Value = 12346967
$ws.Range("K58").Value = 12346967
$ws.Range("M58").Value = -12346764
$ws.Range("H62").Value = 2597.4
$ws.Range("I62").Value = 2177.8
$ws.Range("J62").Value = 3017
$ws.Range("K62").Value = 2177.8
$ws.Range("L62").Value = 3017
$ws.Range("M62").Value = -1553.8
$ws.Range("N62").Value = -4265
$ws.Range("H65").Value = 2597.4
$ws.Range("I65").Value = 2177.8
$ws.Range("J65").Value = 3017
$ws.Range("K65").Value = 10889
$ws.Range("L65").Value = 15085
$ws.Range("M65").Value = -7769
$ws.Range("N65").Value = -21325
$ws.Range("H107").Value = 725.8
$ws.Range("I107").Value = 688.5
$ws.Range("J107").Value = 768.4286
$ws.Range("K107").Value = 688.5
$ws.Range("L107").Value = 768.4286
$ws.Range("M107").Value = 1231.5
$ws.Range("N107").Value = -4608.4286
$ws.Range("H113").Value = 1946.4
$ws.Range("I113").Value = 1722.6
$ws.Range("K113").Value = 1722.6
$ws.Range("M113").Value = 447.4000000000001
$ws.Range("H136").Value = 9806092
$ws.Range("I136").Value = 12346967
$ws.Range("K136").Value = 37040901
$ws.Range("M136").Value = -37038351

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 4003.5
$ws.Range("J88").Value = 4333.3335
$ws.Range("L88").Value = 13000.0005
$ws.Range("N88").Value = -13856.0005
$ws.Range("H91").Value = 4003.5
$ws.Range("J91").Value = 4333.3335
$ws.Range("L91").Value = 13000.0005
$ws.Range("N91").Value = -15964.0005
$ws.Range("H131").Value = 1609.875
$ws.Range("J131").Value = 1662.027
$ws.Range("L131").Value = 4986.081
$ws.Range("N131").Value = -15066.081

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = -4888
$ws.Range("H70").Value = 7831.4443
$ws.Range("I70").Value = 7714.1333
$ws.Range("K70").Value = 7714.1333
$ws.Range("M70").Value = -7444.1333
$ws.Range("H73").Value = 7831.4443
$ws.Range("I73").Value = 7714.1333
$ws.Range("K73").Value = 7714.1333
$ws.Range("M73").Value = -6778.1333
$ws.Range("H132").Value = 27135.568
$ws.Range("I132").Value = 27260.064
$ws.Range("K132").Value = 81780.192
$ws.Range("M132").Value = -79250.192

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3834.1428
$ws.Range("I136").Value = 2383.353
$ws.Range("K136").Value = 7150.059
$ws.Range("M136").Value = -4600.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 113216.73
$ws.Range("J74").Value = 120453.78
$ws.Range("L74").Value = 120453.78
$ws.Range("N74").Value = -122325.78
$ws.Range("H77").Value = 113216.73
$ws.Range("J77").Value = 120453.78
$ws.Range("L77").Value = 361361.34
$ws.Range("N77").Value = -370721.34

Write-Output "Applied all scheduled-runner updates."